# =============================================================
# Expand the KDD17/Stocknet/... distance matrix from 7x7 to 9x9
# (adds SLSYelp + STSGold rows/cols) and start a second small
# "PAD" table below it (rows 12-21) -- matches commit
# "started to add more distances".
# =============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Make sure the text-looking numeric value in B13 is stored as TEXT ---
# (mirrors the existing "diagonal" cells B2, C3, B4, ... which use the
# text number format so values like "-121.3429" do not turn into numbers)
$ws.Range("B13").NumberFormat = "@"

# --- Write every cell value ---
$ws.Range("A1").Value = "PAD"
$ws.Range("B1").Value = "KDD17"
$ws.Range("C1").Value = "Stocknet"
$ws.Range("D1").Value = "EconomyNews"
$ws.Range("E1").Value = "Phrasebank"
$ws.Range("F1").Value = "BBCSport"
$ws.Range("G1").Value = "SLSAmazon"
$ws.Range("H1").Value = "SLSIMBD"
$ws.Range("I1").Value = "SLSYelp"
$ws.Range("J1").Value = "STSGold"
$ws.Range("A2").Value = "KDD17"
$ws.Range("B2").Value = "-0.8304 (±0.01109)"
$ws.Range("C2").Value = "0.5793 (±0.02103)"
$ws.Range("D2").Value = "1.9173 (±0.00757)"
$ws.Range("E2").Value = "1.8697 (±0.00956)"
$ws.Range("F2").Value = "1.9182 (±0.00754)"
$ws.Range("G2").Value = "1.9646 (±0.00397)"
$ws.Range("H2").Value = "1.9658 (±0.00492)"
$ws.Range("I2").Value = "1.9736 (±0.00372)"
$ws.Range("J2").Value = "1.9604 (±0.00508)"
$ws.Range("A3").Value = "Stocknet"
$ws.Range("B3").Value = "0.5793 (±0.02103)"
$ws.Range("C3").Value = "-1.0409 (±0.04034)"
$ws.Range("D3").Value = "1.7979 (±0.01644)"
$ws.Range("E3").Value = "1.7443 (±0.01189)"
$ws.Range("F3").Value = "1.8839 (±0.01401)"
$ws.Range("G3").Value = "1.9295 (±0.01568)"
$ws.Range("H3").Value = "1.9337 (±0.01045)"
$ws.Range("I3").Value = "1.9443 (±0.00845)"
$ws.Range("J3").Value = "1.9460 (±0.00672)"
$ws.Range("A4").Value = "EconomyNews"
$ws.Range("B4").Value = "1.9173 (±0.00757)"
$ws.Range("C4").Value = "1.7979 (±0.01644)"
$ws.Range("D4").Value = "-1.6018 (±0.08220)"
$ws.Range("E4").Value = "1.7978 (±0.03142)"
$ws.Range("F4").Value = "1.8871 (±0.01388)"
$ws.Range("G4").Value = "1.8914 (±0.02574)"
$ws.Range("H4").Value = "1.8965 (±0.03098)"
$ws.Range("I4").Value = "1.9105 (±0.03012)"
$ws.Range("J4").Value = "1.9292 (±0.02288)"
$ws.Range("A5").Value = "Phrasebank"
$ws.Range("B5").Value = "1.8697 (±0.00956)"
$ws.Range("C5").Value = "1.7443 (±0.01189)"
$ws.Range("D5").Value = "1.7978 (±0.03142)"
$ws.Range("E5").Value = "-1.4674 (±0.019566)"
$ws.Range("F5").Value = "1.7842 (±0.01275)"
$ws.Range("G5").Value = "1.8844 (±0.00978)"
$ws.Range("H5").Value = "1.8998 (±0.01852)"
$ws.Range("I5").Value = "1.8998 (±0.01094)"
$ws.Range("J5").Value = "1.9299 (±0.00649)"
$ws.Range("A6").Value = "BBCSport"
$ws.Range("B6").Value = "1.9182 (±0.00754)"
$ws.Range("C6").Value = "1.8839 (±0.01401)"
$ws.Range("D6").Value = "1.8871 (±0.01388)"
$ws.Range("E6").Value = "1.7842 (±0.01275)"
$ws.Range("F6").Value = "-1.4021 (±0.02795)"
$ws.Range("G6").Value = "1.7804 (±0.03478)"
$ws.Range("H6").Value = "1.7451 (±0.02679)"
$ws.Range("I6").Value = "1.7608 (±0.01885)"
$ws.Range("J6").Value = "1.6766 (±0.02538)"
$ws.Range("A7").Value = "SLSAmazon"
$ws.Range("B7").Value = "1.9646 (±0.00397)"
$ws.Range("C7").Value = "1.9295 (±0.01568)"
$ws.Range("D7").Value = "1.8914 (±0.02574)"
$ws.Range("E7").Value = "1.8844 (±0.00978)"
$ws.Range("F7").Value = "1.7804 (±0.03478)"
$ws.Range("G7").Value = "-1.5290 (±0.06082)"
$ws.Range("H7").Value = "1.3460 (±0.07877)"
$ws.Range("I7").Value = "1.2920 (±0.10891)"
$ws.Range("J7").Value = "1.6633 (±0.05061)"
$ws.Range("A8").Value = "SLSIMBD"
$ws.Range("B8").Value = "1.9658 (±0.00492)"
$ws.Range("C8").Value = "1.9337(±0.01045)"
$ws.Range("D8").Value = "1.8965 (±0.03098)"
$ws.Range("E8").Value = "1.8998 (±0.01852)"
$ws.Range("F8").Value = "1.7451 (±0.02679)"
$ws.Range("G8").Value = "1.3460 (±0.07877)"
$ws.Range("H8").Value = "-1.5820 (±0.04367)"
$ws.Range("I8").Value = "1.4260 (±0.06204)"
$ws.Range("J8").Value = "1.6850 (±0.03915)"
$ws.Range("A9").Value = "SLSYelp"
$ws.Range("B9").Value = "1.9736 (±0.00372)"
$ws.Range("C9").Value = "1.9443 (±0.00845)"
$ws.Range("D9").Value = "1.9105 (±0.03012)"
$ws.Range("E9").Value = "1.8998 (±0.01094)"
$ws.Range("F9").Value = "1.7608 (±0.01885)"
$ws.Range("G9").Value = "1.2920 (±0.10891)"
$ws.Range("H9").Value = "1.4260 (±0.06204)"
$ws.Range("I9").Value = "-1.5870 (±0.05438)"
$ws.Range("J9").Value = "1.5941 (±0.02867)"
$ws.Range("A10").Value = "STSGold"
$ws.Range("B10").Value = "1.9604 (±0.00508)"
$ws.Range("C10").Value = "1.9460 (±0.00672)"
$ws.Range("D10").Value = "1.9292 (±0.02288)"
$ws.Range("E10").Value = "1.9299 (±0.00649)"
$ws.Range("F10").Value = "1.6766 (±0.02538)"
$ws.Range("G10").Value = "1.6633 (±0.05061)"
$ws.Range("H10").Value = "1.6850 (±0.03915)"
$ws.Range("I10").Value = "1.5941 (±0.02867)"
$ws.Range("J10").Value = "-1.5553 (±0.02868)"
$ws.Range("B12").Value = "KDD17"
$ws.Range("C12").Value = "Stocknet"
$ws.Range("D12").Value = "EconomyNews"
$ws.Range("E12").Value = "Phrasebank"
$ws.Range("F12").Value = "BBCSport"
$ws.Range("G12").Value = "SLSAmazon"
$ws.Range("H12").Value = "SLSIMBD"
$ws.Range("I12").Value = "SLSYelp"
$ws.Range("J12").Value = "STSGold"
$ws.Range("A13").Value = "KDD17"
$ws.Range("B13").Value = "-121.3429"
$ws.Range("A14").Value = "Stocknet"
$ws.Range("A15").Value = "EconomyNews"
$ws.Range("A16").Value = "Phrasebank"
$ws.Range("A17").Value = "BBCSport"
$ws.Range("A18").Value = "SLSAmazon"
$ws.Range("A19").Value = "SLSIMBD"
$ws.Range("A20").Value = "SLSYelp"
$ws.Range("A21").Value = "STSGold"

# --- Re-apply the workbook's existing named look to every cell we touched, ---
# --- by copying the format from a cell that already has the right style.  ---

# "value" cells (plain distance numbers, centered) -> copy from C2
$ws.Range("C2").Copy()
foreach ($cell in $ws.Range("C2,E2,F2,G2,H2,I2,J2,B3,D3,E3,F3,G3,H3,I3,J3,C4,E4,F4,G4,H4,I4,J4,B5,C5,D5,F5,G5,H5,I5,J5,B6,C6,D6,E6,G6,H6,I6,J6,B7,C7,D7,E7,F7,H7,I7,J7,B8,C8,D8,E8,F8,G8,H8,I8,J8,B9,C9,D9,E9,F9,G9,H9,I9,J9,B10,C10,D10,E10,F10,G10,H10,I10,J10")) { $cell.PasteSpecial(-4122) }

# "label" cells (dataset names, bold + centered) -> copy from B1
$ws.Range("B1").Copy()
foreach ($cell in $ws.Range("B1,C1,D1,E1,F1,G1,H1,I1,J1,A2,A3,A4,A5,A6,A7,A8,A9,A10,B12,C12,D12,E12,F12,G12,H12,I12,J12,A13,A14,A15,A16,A17,A18,A19,A20,A21")) { $cell.PasteSpecial(-4122) }

# "diagonal / text" cells (text number format, centered) -> copy from B2
$ws.Range("B2").Copy()
foreach ($cell in $ws.Range("B2,D2,C3,B4,D4,E5,F6,G7,B13,D13,C14,B15,D15,E16,F17,G18")) { $cell.PasteSpecial(-4122) }

$excel.CutCopyMode = 0

# --- Big bold "PAD" title cells (A1 of each table) ---
foreach ($cell in $ws.Range("A1,A12")) {
  $cell.Font.Name = "Calibri"
  $cell.Font.Bold = $true
  $cell.Font.Size = 26
  $cell.HorizontalAlignment = -4108
}

# --- Row heights for the two header rows ---
$ws.Rows(1).RowHeight = 33.75
$ws.Rows(12).RowHeight = 33.75

# --- Column widths (auto-fit feel for the new/resized columns) ---
$ws.Columns("E:E").ColumnWidth = 19.16
$ws.Columns("H:H").ColumnWidth = 18.02
$ws.Columns("I:J").ColumnWidth = 17.31

# --- Selection matches the author's last click position ---
$ws.Range("E19").Select()
